$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add an hour entry for F7 (numbered tile for the clone), which bumps the running totals in column G
$ws.Range("F7").Value = 1

# Update the saved view state: scroll position and active selection
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Application.ActiveWindow.ScrollRow = 2
$ws.Range("F8").Select()
